$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J mirrors the formatting of column I (same style indices used per row),
# but holds the new 2021 data series.
$ws.Range("I4:I14").Copy() | Out-Null
$ws.Range("J4:J14").PasteSpecial(-4122) | Out-Null

# Populate the new 2021 values in column J
$ws.Range("J4").Value = 2021
$ws.Range("J5").Value = 1.5
$ws.Range("J6").Value = 0.3
$ws.Range("J7").Value = 0.8
$ws.Range("J8").Value = 0.6
$ws.Range("J9").Value = 1.8
$ws.Range("J10").Value = 0.5
$ws.Range("J11").Value = 0.8
$ws.Range("J12").Value = 1.9
$ws.Range("J13").Value = 4.4000000000000004
$ws.Range("J14").Value = 0.4

# Match the post-edit selection state
$ws.Range("L10").Select() | Out-Null
